$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column (C) for rows 2-11 from
# serial date 45171 (2023-09-02) to 45172 (2023-09-03).
foreach ($row in 2..11) {
    $ws.Cells.Item($row, 3).Value = 45172
}
